$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.151.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.173.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.23%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.09%  '
$ws.Range("E7").Value = '  -5.35%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.173.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.726.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.233.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.58%  '
$ws.Range("E16").Value = '  -2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.162.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '416.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("E27").Value = '  -5.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("E35").Value = '  -4.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '155.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.698.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.57%  '
$ws.Range("E39").Value = '  -5.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("E41").Value = '  -3.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("E43").Value = '  -7.02%  '
$ws.Range("E44").Value = '  -5.00%  '
$ws.Range("E45").Value = '  -4.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.11%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '296.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.50%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0262'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("E49").Value = '  -10.09%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0993'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.22%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
